$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "IKHI0000148"
$ws.Range("B1").Value = "CREDIT.ACCT.NO"
$ws.Range("B2").Value = 5000040127

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(2).EntireColumn.AutoFit() | Out-Null

$ws.Range("G17").Select() | Out-Null
